# The commit permutes the 15 data rows (rows 2-16) of the "Artfynd" sheet:
# each row's full record (columns A, B, D, E, F, G, H, I, Q, R, S - the only
# columns that differ row-to-row) moves to a new row position. Columns C, J,
# K..P (minus a couple of genuinely blank cells) and T onward are identical
# across all data rows, so no visible change occurs there.
#
# Mapping of new row -> old row (i.e. new row's data is taken from old row):
#   2<-7  3<-2  4<-15 5<-10 6<-11 7<-3  8<-14 9<-12 10<-8
#   11<-16 12<-6 13<-4 14<-5 15<-13 16<-9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per data row.
$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "Q", "R", "S")

# Columns that hold plain numbers.
$numericCols = @("A", "B", "E", "Q", "R", "S")

# Columns that are naturally text and never look like a bare number (so a
# plain string assignment is stored as text with no extra styling needed).
$plainTextCols = @("D", "F", "G", "H")

# Columns whose text values can look like numbers (e.g. "1", "10") and so
# need special handling to avoid Excel auto-converting them to numbers.
$numericLookingTextCols = @("I")

$firstRow = 2
$lastRow = 16

# New row (key) -> old row (value) it should copy its data from.
$mapping = @{
    2  = 7
    3  = 2
    4  = 15
    5  = 10
    6  = 11
    7  = 3
    8  = 14
    9  = 12
    10 = 8
    11 = 16
    12 = 6
    13 = 4
    14 = 5
    15 = 13
    16 = 9
}

# 1) Snapshot every source cell's value before any writes happen, so that
#    overlapping reads/writes of the permutation don't clobber each other.
$snapshot = @{}
foreach ($r in $firstRow..$lastRow) {
    foreach ($c in $cols) {
        $snapshot["$c$r"] = $ws.Range("$c$r").Value2
    }
}

# 2) Write the snapshotted values back out according to the new<-old mapping.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $val = $snapshot["$c$oldRow"]
        $target = $ws.Range("$c$newRow")

        if ($numericCols -contains $c) {
            $target.Value = $val
        }
        elseif ($numericLookingTextCols -contains $c) {
            # These strings look like plain numbers ("1", "10", ...).
            # A plain assignment would get auto-coerced to a real number by
            # Excel, so prefix with an apostrophe (the normal Excel way of
            # forcing text entry). That marks the cell with a quote-prefix
            # style, so restore the cell's style to "Normal" afterwards to
            # leave no visible formatting residue.
            $target.Value = "'" + $val
            $target.Style = "Normal"
        }
        else {
            # Plain text columns (D,F,G,H) - never look like bare numbers,
            # so a direct assignment is already stored as text.
            $target.Value = $val
        }
    }
}

Write-Host "Row permutation applied."
